$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 126, shifting the existing row 126..378 down to 127..379
$ws.Rows.Item(126).Insert()

# Populate the newly inserted row 126 with a new daily price record
$ws.Cells.Item(126, 1).Value = 3
$ws.Cells.Item(126, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(126, 3).Value = "Coquimbo"
$ws.Cells.Item(126, 4).Value = 44883
$ws.Cells.Item(126, 5).Value = 5
$ws.Cells.Item(126, 6).Value = 100112039
$ws.Cells.Item(126, 7).Value = "Ciboulette"
$ws.Cells.Item(126, 8).Value = "Sin especificar"
$ws.Cells.Item(126, 9).Value = "Primera"
$ws.Cells.Item(126, 10).Value = 120
$ws.Cells.Item(126, 11).Value = 1500
$ws.Cells.Item(126, 12).Value = 1500
$ws.Cells.Item(126, 13).Value = 1500
$ws.Cells.Item(126, 14).Value = "$/docena de atados"
$ws.Cells.Item(126, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(126, 16).Value = 500
$ws.Cells.Item(126, 17).Value = 3
$ws.Cells.Item(126, 18).Value = "Hortaliza"
